$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "260.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.56%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.92%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.710"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.88%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06081"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.14%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.673"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.02%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8461"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.43%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9233"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.59%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.92%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "17.93%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07101"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.44%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03129"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.51%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09073"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.33%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001536"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.78%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006102"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.06%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006134"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.39%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.452"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.148"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.71%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.167"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.24%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.47%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1305"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.86%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.094"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.95%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04235"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.48%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.05%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-9.11%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.14%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001576"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.46%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03872"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.79%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1113"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.38%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004094"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-34.86%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "21.63%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002230"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.42%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005320"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.42%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.14%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05452"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "5.10%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1358"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-46.12%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.14%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.14%"
